# fix bug trong file convertJsonToExcel
# Sheet "Đơn thu nợ": update row 2 (O2, S2), turn the old "Tổng" row 3
# into a new "TN" data row, and append a new "Tổng" totals row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Duplicate the current "Tổng" row (row 3) down to row 4: this keeps the
# blank text-columns (D,E,F,G,H,I,J,L,P,Q) present on the totals row, just
# like the original sheet, while row 3 becomes free to hold the new "TN"
# entry below.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(3).Insert(-4121)   # xlShiftDown

# --- Row 2: existing "TN" entry gets higher payment / discount ---
$ws.Cells.Item(2, 15).Value = 12000000   # O2 Đã thanh toán
$ws.Cells.Item(2, 19).Value = 240000     # S2 Chiết khấu sale chính

# --- Row 3: becomes a new "TN" data row (was the old "Tổng" summary) ---
$ws.Cells.Item(3, 1).Value  = "TN"                  # A3 Tiền tố
$ws.Cells.Item(3, 2).Value  = 197                   # B3 Mã đơn thu nợ
$ws.Cells.Item(3, 3).Value  = 1000000                # C3 Lượng thu
$ws.Cells.Item(3, 4).Value  = "'08-24-2024"         # D3 Ngày thu (force text, not a date)
$ws.Cells.Item(3, 5).Value  = "CẦN THƠ"             # E3 Cơ sở
$ws.Cells.Item(3, 6).Value  = "HD-LUXURY-538"       # F3 Đơn nợ
$ws.Cells.Item(3, 7).Value  = "Nâng mũi"            # G3 Tên dịch vụ
$ws.Cells.Item(3, 8).Value  = "Ngô Xuân Nhi"        # H3 Khách hàng
$ws.Cells.Item(3, 9).Value  = "Cá nhân"             # I3 Nguồn khách
$ws.Cells.Item(3, 10).Value = "Lâm Hoàng Phú"       # J3 Sale chính
$ws.Cells.Item(3, 11).Value = 10000000              # K3 Đơn giá gốc
$ws.Cells.Item(3, 12).Value = "Đỗ Thị Huyền Trân"   # L3 Sale phụ
$ws.Cells.Item(3, 13).Value = 8000000               # M3 Upsale
$ws.Cells.Item(3, 14).Value = 18000000              # N3 Đơn giá
$ws.Cells.Item(3, 15).Value = 12000000              # O3 Đã thanh toán
$ws.Cells.Item(3, 16).Value = "Lâm Thị Mỹ Hằng"     # P3 Bác sĩ 1
$ws.Cells.Item(3, 17).Value = 0                     # Q3 Bác sĩ 2
$ws.Cells.Item(3, 18).Value = 0.1                   # R3 Tỉ lệ chiết khấu sale chính
$ws.Cells.Item(3, 19).Value = 160000                # S3 Chiết khấu sale chính
$ws.Cells.Item(3, 20).Value = 0                     # T3 Tỉ lệ chiết khấu sale phụ
$ws.Cells.Item(3, 21).Value = 0                     # U3 Chiết khấu sale phụ
$ws.Cells.Item(3, 22).Value = 0                     # V3 Tỉ lệ chiết khấu bác sĩ 1
$ws.Cells.Item(3, 23).Value = 0                     # W3 Chiết khấu bác sĩ 1
$ws.Cells.Item(3, 24).Value = 0                     # X3 Tỉ lệ chiết khấu bác sĩ 2
$ws.Cells.Item(3, 25).Value = 0                     # Y3 Chiết khấu bác sĩ 2

# --- Row 4: new "Tổng" totals row (sums rows 2 + 3) ---
# A4, D4-J4, L4, P4, Q4 are already correct (carried over from the copy
# above), only the numeric totals need updating.
$ws.Cells.Item(4, 2).Value  = 2        # B4 Mã đơn thu nợ (count)
$ws.Cells.Item(4, 3).Value  = 2500000  # C4 Lượng thu
$ws.Cells.Item(4, 11).Value = 20000000 # K4 Đơn giá gốc
$ws.Cells.Item(4, 13).Value = 16000000 # M4 Upsale
$ws.Cells.Item(4, 14).Value = 36000000 # N4 Đơn giá
$ws.Cells.Item(4, 15).Value = 24000000 # O4 Đã thanh toán
$ws.Cells.Item(4, 18).Value = 0        # R4 Tỉ lệ chiết khấu sale chính
$ws.Cells.Item(4, 19).Value = 400000   # S4 Chiết khấu sale chính
$ws.Cells.Item(4, 20).Value = 0        # T4 Tỉ lệ chiết khấu sale phụ
$ws.Cells.Item(4, 21).Value = 0        # U4 Chiết khấu sale phụ
$ws.Cells.Item(4, 22).Value = 0        # V4 Tỉ lệ chiết khấu bác sĩ 1
$ws.Cells.Item(4, 23).Value = 0        # W4 Chiết khấu bác sĩ 1
$ws.Cells.Item(4, 24).Value = 0        # X4 Tỉ lệ chiết khấu bác sĩ 2
$ws.Cells.Item(4, 25).Value = 0        # Y4 Chiết khấu bác sĩ 2

# --- Sheet "Lương": propagate the updated discount totals ---
$wsL = $wb.Worksheets.Item(5)
$wsL.Cells.Item(10, 2).Value = 400000                 # B10 Chiết khấu thu nợ tại CẦN THƠ
$wsL.Cells.Item(34, 2).Value = 489285.7142857141      # B34 Tổng lương tại CẦN THƠ
$wsL.Cells.Item(37, 2).Value = 589285.7142857141      # B37 Tổng lương tại HỆ THỐNG
